# Saldo.xlsx update ("Add files via upload"):
#   - Remove account 005277762 / NIVALDO / 56945.46
#   - Remove account 004261201 / ANA / 10116.07
#   - Move account 001761119 / BLUEMETRIX / 136777.36 up, right after the header
#   - Update account 005697554 / MARIA balance from 663423.88 to 86423.88
#
# After the two deletions, row 2 (MARIA) and row 3 (BLUEMETRIX) simply swap
# places, with MARIA's balance becoming the new 86423.88 figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two removed rows bottom-up so earlier row numbers stay valid.
$ws.Rows.Item(7).Delete()   # 004261201 / ANA / 10116.07
$ws.Rows.Item(4).Delete()   # 005277762 / NIVALDO / 56945.46

# Row 2 becomes BLUEMETRIX (previously row 3's contents).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "001761119"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "BLUEMETRIX"
$ws.Range("C2").Value = 136777.36

# Row 3 becomes MARIA with the updated balance.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "005697554"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "MARIA"
$ws.Range("C3").Value = 86423.88
